$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Assign the found Range's .Text directly (rather than using
        # Find.Execute's own Replace:=wdReplaceAll) so Word's AutoCorrect
        # "smart quotes" substitution isn't applied to the new text.
        $rng.Text = $new
    } else {
        Write-Output "NOT FOUND: $old"
    }
    return $found
}

Replace-Text "You'll receive a daily notification to remind you to complete your lesson. And if you miss it, it is also okay! You can always return to ParentText anytime to catch up on your lesson." "Jy sal 'n daaglikse kennisgewing ontvang om jou te herinner om jou les te voltooi. En as jy dit mis, moenie bekommerd wees nie! Jy kan enige tyd na ParentText terugkeer om jou les in te haal."

Replace-Text "Each lesson is a mix of quizzes, comics, tips, and a fun activity to try at home with your child or family." "Elke les is 'n mengsel van toetse, strokiesprente, wenke en 'n lekker aktiwiteit om by die huis saam met jou kind of gesin te probeer."

Replace-Text "If you are ever stuck or need help, type MENU or HELP at the end of your lessons to get more support. " "As jy ooit vas hak of hulp nodig het, klik op MENU of HELP aan die einde van jou lesse om meer support te kry. "

Replace-Text "When you type HELP anytime, you can get information about resources in your community to address family violence, sexual violence, mental health, or other emergencies. " "Jy kan enige tyd HELP klik en inligting kry oor hulpbronne in jou gemeenskap om gesinsgeweld, seksuele geweld, geestesgesondheid of ander noodgevalle aan te spreek. "

Replace-Text "Your information here is safe: Nothing will be shared without your permission and will not be sold for profit. The messages you send are encrypted and locked in a secure server. " "Jou inligting is safe hier: Niks sal gedeel word sonder jou toestemming nie en sal nie vir wins verkoop word nie. Die boodskappe wat jy stuur, is geïnkripteer en gesluit in 'n veilige bediener. "

Replace-Text "Remember, anyone with access to your unlocked phone can view your messages. So, if you send sensitive information and are worried, delete the messages from your phone. " "Onthou, enige iemand met toegang tot jou ongeslote foon, kan na jou boodskappe kyk. So, as jy sensitiewe inligting stuur en is jy worried, verwyder die boodskappe van jou foon. "

Replace-Text "Being here shows how much you care about providing the best support for your child. " "Deur hier te wees, wys hoeveel jy omgee om die beste support vir your child gee. "

Replace-Text "It is what you do with your child that will really make a difference. " "Dit is wat jy doen saam met your child wat regtig 'n verskil maak. "

Replace-Text "ParentText will provide tips through lessons to help you with your relationship with your child. It is up to you to put these tips into practice!" "ParentText sal regdeur die lesse wenke gee to help jou met jou relationship met your child. It is up to you to put these tips into practice!"
